# Output Currency Conversion Factors.xlsx -- update to CPL's work-to-date
# (U.S. model, 2020-dollar basis) per commit "Initial update with CPL's work-to-date"

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

# --- 1. Re-label the "Notes on the Brazilian adaptation" block as the
#        U.S.-model heading, matching the plain-bold style already used
#        by the other section headers (A1/A5/A7), then drop the old
#        "For the Brazil model:" line (row 16) entirely so everything
#        below it shifts up by one row.
$about.Range("A1").Copy()
$about.Range("A15").PasteSpecial(-4122)   # xlPasteFormats
$about.Range("A15").Value = "For the U.S. model:"
$about.Rows("16:16").Delete()

# After the delete: old row17->16 (blank spacer), old18->17 (Large Output
# Currency Unit label), old19->18 (billion ... dollars), old20->19 (blank),
# old21->20 (Medium label), old22->21 (million ... dollars), old23->22
# (blank), old24->23 (Small label), old25->24 (... dollars), old26->25
# (blank), old27->26 (conversion value + note).

# --- 2. Update the currency-year labels from 2019 to 2020.
$about.Range("A18").Value = "billion 2020 dollars"
$about.Range("A21").Value = "million 2020 dollars"
$about.Range("A24").Value = "2020 dollars"

# --- 3. Update the conversion factor itself.
$about.Range("A26").Value = 0.88711067149387013

# --- 4. Append the new explanatory note below the factor (rows 28-33).
$about.Range("B28").Value = 'Recall, this variable is "dollars per large/medium/small currency output unit"'
$about.Range("B29").Value = 'which in this case is "2012 dollars per 2020 dollar."'
$about.Range("B30").Value = "2012 dollars are worth more than 2020 dollars, so we need a"
$about.Range("B31").Value = "value less than 1 in this variable."
$about.Range("B32").Value = "This is why we multiply, not divide, by the conversion"
$about.Range("B33").Value = "factor above."

# --- 5. Update the selection/active cell to match the author's final
#        cursor position. The other three sheets go back to showing the
#        default top-left cell (no special selection persisted for them).
$locu = $wb.Worksheets.Item("OCCF-DpLOCU")
$locu.Range("A1").Select()

$about.Activate()
$about.Range("B31").Select()

# The OCCF-DpLOCU / OCCF-DpMOCU / OCCF-DpSOCU sheets reference
# About!$A$27 (now About!$A$26 after the row delete above -- Excel keeps
# that reference correct automatically) and their formulas/cached values
# recalc on their own, so no further edits are required there.
